# [PV-350][WIP] Replace hard coding of visual height with calculated value
#
# The "Row ID" column header on the PV-Test-03 plan sheet is renamed to "Id".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-03")

$ws.Range("A1").Value = "Id"
